# Add new columns I (I0) and J (IF) to Sheet1, filling header + data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy formatting (bold font, borders, centered alignment) from the existing
# "IP" header cell (H1) onto the two new header cells, then overwrite the text.
$ws.Cells.Item(1, 8).Copy($ws.Cells.Item(1, 9))
$ws.Cells.Item(1, 9).Value = "I0"

$ws.Cells.Item(1, 8).Copy($ws.Cells.Item(1, 10))
$ws.Cells.Item(1, 10).Value = "IF"

# --- Data rows (rows 2..47), values for columns I and J ---
$rows = @(
    @{Row=2; I=1; J=1},
    @{Row=3; I=7; J=7},
    @{Row=4; I=7; J=7},
    @{Row=5; I=8; J=8},
    @{Row=6; I=7; J=8},
    @{Row=7; I=7; J=7},
    @{Row=8; I=6; J=6},
    @{Row=9; I=11; J=13},
    @{Row=10; I=7; J=7},
    @{Row=11; I=7; J=9},
    @{Row=12; I=7; J=7},
    @{Row=13; I=1; J=4},
    @{Row=14; I=7; J=7},
    @{Row=15; I=8; J=8},
    @{Row=16; I=1; J=4},
    @{Row=17; I=9; J=9},
    @{Row=18; I=9; J=9},
    @{Row=19; I=6; J=6},
    @{Row=20; I=1; J=6},
    @{Row=21; I=1; J=5},
    @{Row=22; I=3; J=7},
    @{Row=23; I=3; J=6},
    @{Row=24; I=1; J=6},
    @{Row=25; I=1; J=4},
    @{Row=26; I=1; J=4},
    @{Row=27; I=1; J=4},
    @{Row=28; I=1; J=5},
    @{Row=29; I=1; J=6},
    @{Row=30; I=1; J=6},
    @{Row=31; I=1; J=6},
    @{Row=32; I=1; J=5},
    @{Row=33; I=1; J=6},
    @{Row=34; I=1; J=6},
    @{Row=35; I=1; J=6},
    @{Row=36; I=1; J=4},
    @{Row=37; I=1; J=7},
    @{Row=38; I=1; J=5},
    @{Row=39; I=1; J=8},
    @{Row=40; I=1; J=6},
    @{Row=41; I=1; J=5},
    @{Row=42; I=1; J=5},
    @{Row=43; I=1; J=4},
    @{Row=44; I=1; J=5},
    @{Row=45; I=1; J=4},
    @{Row=46; I=1; J=3},
    @{Row=47; I=1; J=2}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Cells.Item($r.Row, 10).Value = $r.J
}

# Dimension will update automatically to reflect the new used range (A1:J47).
